{"js": "// Apply the fiscal-report field updates described by the diff:\n//  - Immeuble ID: 4 -> 1\n//  - Propri\u00e9taire name \"Luigi Brothers\" cleared (stray trailing space kept)\n//  - Contact line: email \"ezfthgfrdezd\" added, phone number cleared\n//  - Revenus locatifs totaux : 8300.0 -> 750.0\n//  - D\u00e9penses totales (charges et travaux) : 108.4 -> 253.4\n//  - Taxes et factures impay\u00e9es : 2400.0 -> 0.0\n\nconst body = context.document.body;\n\n// Each [needle, replacement] pair is a literal, case-sensitive, unique\n// substring of the document body text (verified against before.docx).\nconst replacements = [\n  [\n    \"Rapport de D\u00e9claration Fiscale pour l'Immeuble ID: 4\",\n    \"Rapport de D\u00e9claration Fiscale pour l'Immeuble ID: 1\"\n  ],\n  [\n    \"Propri\u00e9taire : Luigi Brothers\",\n    \"Propri\u00e9taire :  \"\n  ],\n  [\n    \"Contact : , 07 58 47 61 25\",\n    \"Contact : ezfthgfrdezd, \"\n  ],\n  [\n    \"Revenus locatifs totaux : 8300.0 \u20ac\",\n    \"Revenus locatifs totaux : 750.0 \u20ac\"\n  ],\n  [\n    \"D\u00e9penses totales (charges et travaux) : 108.4 \u20ac\",\n    \"D\u00e9penses totales (charges et travaux) : 253.4 \u20ac\"\n  ],\n  [\n    \"Taxes et factures impay\u00e9es : 2400.0 \u20ac\",\n    \"Taxes et factures impay\u00e9es : 0.0 \u20ac\"\n  ]\n];\n\nfor (const [needle, replacement] of replacements) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${needle}`);\n  }\n\n  // Replace every occurrence (expected to be exactly one for each needle).\n  for (const hit of results.items) {\n    hit.insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the fiscal-report field updates described by the diff:\n#  - Immeuble ID: 4 -> 1\n#  - Propri\u00e9taire name \"Luigi Brothers\" cleared (stray trailing space kept)\n#  - Contact line: email \"ezfthgfrdezd\" added, phone number cleared\n#  - Revenus locatifs totaux : 8300.0 -> 750.0\n#  - D\u00e9penses totales (charges et travaux) : 108.4 -> 253.4\n#  - Taxes et factures impay\u00e9es : 2400.0 -> 0.0\n\n$d = $word.ActiveDocument\n\n# Avoid Word's \"smart quotes\" AutoFormat/AutoCorrect kicking in when text is\n# written back (it would turn the straight apostrophe in \"l'Immeuble\" into a\n# curly one). Belt-and-braces: turn it off, and also write the replacement\n# through Range.Text (not Find.Replacement.Text) which never triggers it.\ntry { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}\ntry { $word.Options.AutoFormatReplaceQuotes = $false } catch {}\n\nfunction Replace-FirstMatch {\n    param(\n        $Document,\n        [string]$Needle,\n        [string]$Replacement\n    )\n\n    $range = $Document.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $Needle\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n\n    $found = $range.Find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $Needle\"\n    }\n\n    # $range now spans exactly the matched text (Find.Execute collapses /\n    # extends the range to the hit) \u2014 overwrite it directly so formatting on\n    # the surrounding run is preserved and no AutoCorrect substitution runs.\n    $range.Text = $Replacement\n}\n\nReplace-FirstMatch $d \"Rapport de D\u00e9claration Fiscale pour l'Immeuble ID: 4\" \"Rapport de D\u00e9claration Fiscale pour l'Immeuble ID: 1\"\nReplace-FirstMatch $d \"Propri\u00e9taire : Luigi Brothers\" \"Propri\u00e9taire :  \"\nReplace-FirstMatch $d \"Contact : , 07 58 47 61 25\" \"Contact : ezfthgfrdezd, \"\nReplace-FirstMatch $d \"Revenus locatifs totaux : 8300.0 \u20ac\" \"Revenus locatifs totaux : 750.0 \u20ac\"\nReplace-FirstMatch $d \"D\u00e9penses totales (charges et travaux) : 108.4 \u20ac\" \"D\u00e9penses totales (charges et travaux) : 253.4 \u20ac\"\nReplace-FirstMatch $d \"Taxes et factures impay\u00e9es : 2400.0 \u20ac\" \"Taxes et factures impay\u00e9es : 0.0 \u20ac\"\n"}
